$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 8412
$ws.Range("F6").Value = 534
$ws.Range("F7").Value = 7379
$ws.Range("F8").Value = 7379
$ws.Range("F9").Value = 593
$ws.Range("F10").Value = 513
$ws.Range("F15").Value = 235
$ws.Range("F18").Value = 157
$ws.Range("F19").Value = 12234
$ws.Range("F22").Value = 2497
$ws.Range("F23").Value = 3580
$ws.Range("F26").Value = 2937
$ws.Range("F27").Value = 113
$ws.Range("F28").Value = 113
$ws.Range("F29").Value = 16
$ws.Range("F30").Value = 37
$ws.Range("F31").Value = 3362
$ws.Range("F32").Value = 80
$ws.Range("F34").Value = 1729
$ws.Range("F35").Value = 81
$ws.Range("F36").Value = 138
$ws.Range("F37").Value = 6044
$ws.Range("F38").Value = 101
$ws.Range("F39").Value = 1842
$ws.Range("F40").Value = 1260
$ws.Range("F41").Value = 35
$ws.Range("F42").Value = 917
$ws.Range("F47").Value = 1126
$ws.Range("F49").Value = 1592
$ws.Range("F50").Value = 26

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 257
$ws.Range("F26").Value = 7

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 329
$ws.Range("F3").Value = 479

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 329
$ws.Range("F7").Value = 8413
$ws.Range("F10").Value = 534
$ws.Range("F11").Value = 7379
$ws.Range("F12").Value = 593
$ws.Range("F13").Value = 513
$ws.Range("F16").Value = 235
$ws.Range("F17").Value = 257
$ws.Range("F19").Value = 157
$ws.Range("F21").Value = 12235
$ws.Range("F25").Value = 2497
$ws.Range("F26").Value = 2497
$ws.Range("F27").Value = 3580
$ws.Range("F28").Value = 113
$ws.Range("F29").Value = 113
$ws.Range("F30").Value = 16
$ws.Range("F31").Value = 37
$ws.Range("F33").Value = 3362
$ws.Range("F35").Value = 1729
$ws.Range("F36").Value = 81
$ws.Range("F37").Value = 138
$ws.Range("F38").Value = 6044
$ws.Range("F40").Value = 101
$ws.Range("F41").Value = 1842
$ws.Range("F43").Value = 1260
$ws.Range("F44").Value = 35
$ws.Range("F45").Value = 917
$ws.Range("F48").Value = 1126
$ws.Range("F50").Value = 1592
$ws.Range("F51").Value = 26

$wb.Save()